$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 11 (openai/gpt-oss-120b)
$ws.Range("I11").Value = 1
$ws.Range("J11").Value = 0.001
$ws.Range("K11").Value = 508
$ws.Range("L11").Value = 0.00254

# Row 12 (openai/gpt-oss-20b)
$ws.Range("I12").Value = 1
$ws.Range("J12").Value = 0.001
$ws.Range("K12").Value = 508
$ws.Range("L12").Value = 0.00254
